# Disaggregation of commodity Copper
#
# 1) Rename the shared text "Copper ores and concentrates" -> "Copper"
#    (row 7 / column C on every yearly sheet references this string).
# 2) Because the "Copper" commodity row is being disaggregated into the
#    three sector columns, the existing D/E/F (Photovoltaic / Offshore
#    wind / Onshore wind) figures for rows 5 (Neodymium), 7 (Copper) and
#    8 (Raw silicon) are rotated one column to the right:
#       new D = old F
#       new E = old D
#       new F = old E
#    Row 6 (Dysprosium) is always zero across the whole workbook so the
#    rotation is a no-op there, but we still include it for completeness.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- rename the commodity label in column C, row 7 ---
    $ws.Range("C7").Value2 = "Copper"

    # --- rotate the D/E/F figures for rows 5, 6, 7, 8 ---
    foreach ($row in 5, 6, 7, 8) {
        $dCell = $ws.Cells.Item($row, 4)
        $eCell = $ws.Cells.Item($row, 5)
        $fCell = $ws.Cells.Item($row, 6)

        $dVal = $dCell.Value2
        $eVal = $eCell.Value2
        $fVal = $fCell.Value2

        $dCell.Value2 = $fVal
        $eCell.Value2 = $dVal
        $fCell.Value2 = $eVal
    }
}
